# goal-demo-schedule.xlsx edit:
#   - Lecturer name "Johnson" -> "Winchester" (cell E2)
#   - Course code "7C" -> "A8" (cell E3)
#   - Move/leave the active selection on E3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "Winchester"
$ws.Range("E3").Value = "A8"

$ws.Range("E3").Select()
